# The workbook lists NBA players, their positions, and their teams in
# columns A, B, C (row 1 = headers). This edit swaps the data for
# "Jalen Green" (row 2) and "Russell Westbrook" (row 14), i.e. the two
# rows trade places in the list (player name + team swap; position is
# the same "PG,SG" for both so it is unaffected either way).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values of the two rows before overwriting anything.
$row2Name = $ws.Range("A2").Value2
$row2Pos  = $ws.Range("B2").Value2
$row2Team = $ws.Range("C2").Value2

$row14Name = $ws.Range("A14").Value2
$row14Pos  = $ws.Range("B14").Value2
$row14Team = $ws.Range("C14").Value2

# Swap the two rows.
$ws.Range("A2").Value = $row14Name
$ws.Range("B2").Value = $row14Pos
$ws.Range("C2").Value = $row14Team

$ws.Range("A14").Value = $row2Name
$ws.Range("B14").Value = $row2Pos
$ws.Range("C14").Value = $row2Team
